$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the backup code currently in A11 up to A2 (replacing the old A2 value),
# and clear out the now-unused rows A3, A4, A11.
$ws.Range("A2").Value2 = $ws.Range("A11").Value2

$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("A11").ClearContents()

# Update the active selection to A3, matching the saved workbook state.
$ws.Range("A3").Select()
